# Update "想去人数" (want-to-go count) figures for a handful of events
# on both the "展览" and "全部类型" worksheets, matching the latest
# gh-pages data generation run.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 291
    $ws.Range("F4").Value = 10232
    $ws.Range("F18").Value = 357
    $ws.Range("F21").Value = 1589
}
